$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the formatting of H1.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I2:I28 with 1, and J2:J28 with a copy of the corresponding H value.
for ($r = 2; $r -le 28; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}

# Row 29 is the exception: I29 = 5 and J29 = 5 (H29 stays 1).
$ws.Cells.Item(29, 9).Value = 5
$ws.Cells.Item(29, 10).Value = 5

# Extend the sheet's used-range dimension to cover the new columns.
$ws.UsedRange | Out-Null
